$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H110").Value = 60000
$ws.Range("J110").Value = 60000
$ws.Range("L110").Value = 60000
$ws.Range("N110").Value = -68180

$ws.Range("H141").Value = 1122389.1
$ws.Range("I141").Value = 1475670
$ws.Range("J141").Value = 3666.3333
$ws.Range("K141").Value = 4427010
$ws.Range("L141").Value = 10998.9999
$ws.Range("M141").Value = -4421830
$ws.Range("N141").Value = -21358.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3591.4666
$ws.Range("I32").Value = 2721.9395
$ws.Range("K32").Value = 2721.9395
$ws.Range("M32").Value = -2434.9395

$ws.Range("H45").Value = 2986.6956
$ws.Range("I45").Value = 3111.75
$ws.Range("K45").Value = 3111.75
$ws.Range("M45").Value = -2734.75

$ws.Range("H61").Value = 3490.762
$ws.Range("I61").Value = 2499.3845
$ws.Range("K61").Value = 2499.3845
$ws.Range("M61").Value = -2287.3845

$ws.Range("H74").Value = 1984.1818
$ws.Range("I74").Value = 702
$ws.Range("K74").Value = 702
$ws.Range("M74").Value = 172

$ws.Range("H77").Value = 1984.1818
$ws.Range("I77").Value = 702
$ws.Range("K77").Value = 3510
$ws.Range("M77").Value = 858

$ws.Range("H122").Value = 28252.215
$ws.Range("I122").Value = 35511.91
$ws.Range("K122").Value = 106535.73
$ws.Range("M122").Value = -104085.73

$ws.Range("H136").Value = 3490.762
$ws.Range("I136").Value = 2499.3845
$ws.Range("K136").Value = 7498.1535
$ws.Range("M136").Value = -4948.1535

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1389.5454
$ws.Range("I31").Value = 869.7778
$ws.Range("J31").Value = 2013.2667
$ws.Range("K31").Value = 869.7778
$ws.Range("L31").Value = 2013.2667
$ws.Range("M31").Value = -574.7778
$ws.Range("N31").Value = -2603.2667

$ws.Range("H34").Value = 1389.5454
$ws.Range("I34").Value = 869.7778
$ws.Range("J34").Value = 2013.2667
$ws.Range("K34").Value = 869.7778
$ws.Range("L34").Value = 2013.2667
$ws.Range("M34").Value = -667.7778
$ws.Range("N34").Value = -2417.2667

$ws.Range("H105").Value = 1400
$ws.Range("I105").Value = 1566.6666
$ws.Range("K105").Value = 1566.6666
$ws.Range("M105").Value = 180.3334

$ws.Range("H132").Value = 2312.3845
$ws.Range("I132").Value = 1608.7646
$ws.Range("J132").Value = 3641.4443
$ws.Range("K132").Value = 4826.293799999999
$ws.Range("L132").Value = 10924.3329
$ws.Range("M132").Value = -2296.293799999999
$ws.Range("N132").Value = -15984.3329

$ws.Range("H134").Value = 2117.3333
$ws.Range("I134").Value = 1017.8333
$ws.Range("K134").Value = 3053.4999
$ws.Range("M134").Value = -518.4998999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1608.279
$ws.Range("J68").Value = 1911.9688
$ws.Range("L68").Value = 5735.9064
$ws.Range("N68").Value = -7357.9064

$ws.Range("H69").Value = 2498
$ws.Range("I69").Value = 2498
$ws.Range("K69").Value = 7494
$ws.Range("M69").Value = -6683

$ws.Range("H71").Value = 1608.279
$ws.Range("J71").Value = 1911.9688
$ws.Range("L71").Value = 17207.7192
$ws.Range("N71").Value = -25319.7192

$ws.Range("H72").Value = 2498
$ws.Range("I72").Value = 2498
$ws.Range("K72").Value = 22482
$ws.Range("M72").Value = -18426

$ws.Range("H129").Value = 56360.54
$ws.Range("I129").Value = 660
$ws.Range("J129").Value = 121344.5
$ws.Range("K129").Value = 1980
$ws.Range("L129").Value = 364033.5
$ws.Range("M129").Value = 3020
$ws.Range("N129").Value = -374033.5

$ws.Range("H131").Value = 6678440
$ws.Range("J131").Value = 12411.972
$ws.Range("L131").Value = 37235.916
$ws.Range("N131").Value = -47315.916

$ws.Range("H132").Value = 1037.7222
$ws.Range("J132").Value = 1026.4117
$ws.Range("L132").Value = 9237.705300000001
$ws.Range("N132").Value = -14297.7053

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 5614001
$ws.Range("I12").Value = 7000000
$ws.Range("K12").Value = 7000000
$ws.Range("M12").Value = -6999860

$ws.Range("H80").Value = 12001.667
$ws.Range("I80").Value = 9001.286
$ws.Range("J80").Value = 16202.2
$ws.Range("K80").Value = 9001.286
$ws.Range("L80").Value = 16202.2
$ws.Range("M80").Value = -8003.286
$ws.Range("N80").Value = -18198.2

$ws.Range("H83").Value = 12001.667
$ws.Range("I83").Value = 9001.286
$ws.Range("J83").Value = 16202.2
$ws.Range("K83").Value = 45006.43
$ws.Range("L83").Value = 81011
$ws.Range("M83").Value = -40014.43
$ws.Range("N83").Value = -90995

$ws.Range("H132").Value = 856638.5600000001
$ws.Range("I132").Value = 1242074.5
$ws.Range("J132").Value = 3173.4285
$ws.Range("K132").Value = 3726223.5
$ws.Range("L132").Value = 9520.2855
$ws.Range("M132").Value = -3723693.5
$ws.Range("N132").Value = -14580.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7826.091
$ws.Range("I40").Value = 3727.1428
$ws.Range("K40").Value = 3727.1428
$ws.Range("M40").Value = -3591.1428

$ws.Range("H122").Value = 9143.333000000001
$ws.Range("I122").Value = 7815.3
$ws.Range("J122").Value = 11799.4
$ws.Range("K122").Value = 23445.9
$ws.Range("L122").Value = 35398.2
$ws.Range("M122").Value = -20995.9
$ws.Range("N122").Value = -40298.2

$ws.Range("H127").Value = 35570.57
$ws.Range("J127").Value = 35570.57
$ws.Range("L127").Value = 35570.57
$ws.Range("N127").Value = -45490.57

$ws.Range("H132").Value = 2163.4546
$ws.Range("I132").Value = 853
$ws.Range("J132").Value = 2733.2173
$ws.Range("K132").Value = 2559
$ws.Range("L132").Value = 8199.651899999999
$ws.Range("M132").Value = -29
$ws.Range("N132").Value = -13259.6519

$ws.Range("H136").Value = 4778.5713
$ws.Range("I136").Value = 3666
$ws.Range("K136").Value = 10998
$ws.Range("M136").Value = -8448

$ws.Range("H139").Value = 39900
$ws.Range("J139").Value = 39900
$ws.Range("L139").Value = 39900
$ws.Range("N139").Value = -50180

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 271
$ws.Range("I100").Value = 271
$ws.Range("K100").Value = 542
$ws.Range("M100").Value = -1

$ws.Range("H122").Value = 39990.88
$ws.Range("I122").Value = 46871.43
$ws.Range("K122").Value = 140614.29
$ws.Range("M122").Value = -138164.29

$ws.Range("H132").Value = 2501.3704
$ws.Range("I132").Value = 2029.3334
$ws.Range("J132").Value = 3091.4167
$ws.Range("K132").Value = 6088.0002
$ws.Range("L132").Value = 9274.250100000001
$ws.Range("M132").Value = -3558.0002
$ws.Range("N132").Value = -14334.2501

$ws.Range("H136").Value = 26459518
$ws.Range("I136").Value = 61733540
$ws.Range("K136").Value = 185200620
$ws.Range("M136").Value = -185198070
